$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.921.64"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.354.93"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.80"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.14%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.668"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.87%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.26"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -2.24%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.02%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +1.40%  "

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.24%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.03"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.78%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.67"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.92%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.32"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.77%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.705.77"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.57%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "16.44"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.86%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.910"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.356.81"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.64%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.838.31"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.03%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.73%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.72"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.76%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.69"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.05"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -0.37%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +14.95%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.04%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.68%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.64"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.96%  "

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.97%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.71"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.04%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.34"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.35%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.46%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.137"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.59%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.21"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -3.23%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.50"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +2.24%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.58%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -2.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0278"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "68.35"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +27.72%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.111"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +10.13%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.15"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +15.44%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +2.91%  "

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.06"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.50%  "

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.203"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.44%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.44%  "

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "99.80"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.62%  "

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -3.59%  "
